# Apply corrections to EmpleadosConSalarios worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension of data: we now have 5 data rows (rows 2-6) plus header (row 1).
# Row 2: Julian / 19 / 1200 / 36 / 8 / PROGRAMADOR / 43200
$ws.Range("A2").Value = "Julian"
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = 1200
$ws.Range("D2").Value = 36
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "PROGRAMADOR"
$ws.Range("G2").Value = 43200

# Row 3: Julian / 19 / 1200 / 36 / 8 / PROGRAMADOR / 43200
$ws.Range("A3").Value = "Julian"
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 1200
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = "PROGRAMADOR"
$ws.Range("G3").Value = 43200

# Row 4: Juan / 20 / 1200 / 40 / 9 / PROGRAMADOR / 48000
$ws.Range("A4").Value = "Juan"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 1200
$ws.Range("D4").Value = 40
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = "PROGRAMADOR"
$ws.Range("G4").Value = 48000

# Row 5: Juan / 20 / 1200 / 40 / 9 / PROGRAMADOR / 48000
$ws.Range("A5").Value = "Juan"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 1200
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "PROGRAMADOR"
$ws.Range("G5").Value = 48000

# Row 6 (new): Alexis  / 19 / 1000 / 40 / 10 / ANALISTA / 50000
$ws.Range("A6").Value = "Alexis "
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "ANALISTA"
$ws.Range("G6").Value = 50000
